$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 40
$ws.Range("H40").Value = 1216.6666
$ws.Range("I40").Value = 1125
$ws.Range("J40").Value = 1400
$ws.Range("K40").Value = 1125
$ws.Range("L40").Value = 1400
$ws.Range("M40").Value = -950
$ws.Range("N40").Value = -1750

# Row 107
$ws.Range("H107").Value = 1068.2
$ws.Range("I107").Value = 1109.1305
$ws.Range("J107").Value = 597.5
$ws.Range("K107").Value = 1109.1305
$ws.Range("L107").Value = 597.5
$ws.Range("M107").Value = 810.8695
$ws.Range("N107").Value = -4437.5

# Row 116
$ws.Range("H116").Value = 11591.272
$ws.Range("I116").Value = 13611.556
$ws.Range("J116").Value = 2500
$ws.Range("K116").Value = 13611.556
$ws.Range("L116").Value = 2500
$ws.Range("M116").Value = -10169.556
$ws.Range("N116").Value = -9384

# Row 132
$ws.Range("H132").Value = 2729.1777
$ws.Range("I132").Value = 2734.5789
$ws.Range("J132").Value = 2699.8572
$ws.Range("K132").Value = 8203.736699999999
$ws.Range("L132").Value = 8099.571599999999
$ws.Range("M132").Value = -5673.736699999999
$ws.Range("N132").Value = -13159.5716

# Row 137
$ws.Range("H137").Value = 1123.9242
$ws.Range("I137").Value = 952.9804
$ws.Range("J137").Value = 1705.1333
$ws.Range("K137").Value = 2858.9412
$ws.Range("L137").Value = 5115.3999
$ws.Range("M137").Value = -308.9412000000002
$ws.Range("N137").Value = -10215.3999

# Row 138
$ws.Range("H138").Value = 1322.27
$ws.Range("I138").Value = 608.125
$ws.Range("J138").Value = 1981.4807
$ws.Range("K138").Value = 1824.375
$ws.Range("L138").Value = 5944.4421
$ws.Range("M138").Value = 3315.625
$ws.Range("N138").Value = -16224.4421

# Row 141
$ws.Range("H141").Value = 3024.9783
$ws.Range("I141").Value = 1087.6177
$ws.Range("J141").Value = 8514.166999999999
$ws.Range("K141").Value = 3262.8531
$ws.Range("L141").Value = 25542.501
$ws.Range("M141").Value = 1917.1469
$ws.Range("N141").Value = -35902.501


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 1188429.6
$ws.Range("I32").Value = 1395783.2
$ws.Range("J32").Value = 22065.875
$ws.Range("K32").Value = 1395783.2
$ws.Range("L32").Value = 22065.875
$ws.Range("M32").Value = -1395496.2
$ws.Range("N32").Value = -22639.875

# Row 45
$ws.Range("H45").Value = 2044.05
$ws.Range("I45").Value = 2034.7858
$ws.Range("J45").Value = 2065.6667
$ws.Range("K45").Value = 2034.7858
$ws.Range("L45").Value = 2065.6667
$ws.Range("M45").Value = -1657.7858
$ws.Range("N45").Value = -2819.6667

# Row 61
$ws.Range("H61").Value = 2371.5667
$ws.Range("I61").Value = 2087.3333
$ws.Range("J61").Value = 2899.4285
$ws.Range("K61").Value = 2087.3333
$ws.Range("L61").Value = 2899.4285
$ws.Range("M61").Value = -1875.3333
$ws.Range("N61").Value = -3323.4285

# Row 74
$ws.Range("H74").Value = 1079.2954
$ws.Range("I74").Value = 859.75
$ws.Range("J74").Value = 1342.75
$ws.Range("K74").Value = 859.75
$ws.Range("L74").Value = 1342.75
$ws.Range("M74").Value = 14.25
$ws.Range("N74").Value = -3090.75

# Row 77
$ws.Range("H77").Value = 1079.2954
$ws.Range("I77").Value = 859.75
$ws.Range("J77").Value = 1342.75
$ws.Range("K77").Value = 4298.75
$ws.Range("L77").Value = 6713.75
$ws.Range("M77").Value = 69.25
$ws.Range("N77").Value = -15449.75

# Row 97
$ws.Range("H97").Value = 901.9167
$ws.Range("I97").Value = 811.75
$ws.Range("J97").Value = 1352.75
$ws.Range("K97").Value = 811.75
$ws.Range("L97").Value = 1352.75
$ws.Range("M97").Value = -315.75
$ws.Range("N97").Value = -2344.75

# Row 136
$ws.Range("H136").Value = 2371.5667
$ws.Range("I136").Value = 2087.3333
$ws.Range("J136").Value = 2899.4285
$ws.Range("K136").Value = 6261.999899999999
$ws.Range("L136").Value = 8698.2855
$ws.Range("M136").Value = -3711.999899999999
$ws.Range("N136").Value = -13798.2855


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 140
$ws.Range("H140").Value = 57923.332
$ws.Range("J140").Value = 57923.332
$ws.Range("L140").Value = 57923.332
$ws.Range("N140").Value = -68283.33199999999


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 6263.484
$ws.Range("I31").Value = 2137.5557
$ws.Range("J31").Value = 7951.364
$ws.Range("K31").Value = 2137.5557
$ws.Range("L31").Value = 7951.364
$ws.Range("M31").Value = -1842.5557
$ws.Range("N31").Value = -8541.364

# Row 34
$ws.Range("H34").Value = 6263.484
$ws.Range("I34").Value = 2137.5557
$ws.Range("J34").Value = 7951.364
$ws.Range("K34").Value = 2137.5557
$ws.Range("L34").Value = 7951.364
$ws.Range("M34").Value = -1935.5557
$ws.Range("N34").Value = -8355.364

# Row 58
$ws.Range("H58").Value = 1091.3182
$ws.Range("I58").Value = 795.38464
$ws.Range("K58").Value = 795.38464
$ws.Range("M58").Value = -592.38464

# Row 99
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

# Row 122
$ws.Range("H122").Value = 2052.4443
$ws.Range("I122").Value = 2236
$ws.Range("K122").Value = 6708
$ws.Range("M122").Value = -4258

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

# Row 132
$ws.Range("H132").Value = 2977572.5
$ws.Range("I132").Value = 1341.3684
$ws.Range("J132").Value = 9260727
$ws.Range("K132").Value = 4024.1052
$ws.Range("L132").Value = 27782181
$ws.Range("M132").Value = -1494.1052
$ws.Range("N132").Value = -27787241

# Row 134
$ws.Range("H134").Value = 3176.4119
$ws.Range("I134").Value = 3501.9756
$ws.Range("J134").Value = 1841.6
$ws.Range("K134").Value = 10505.9268
$ws.Range("L134").Value = 5524.799999999999
$ws.Range("M134").Value = -7970.926800000001
$ws.Range("N134").Value = -10594.8

# Row 136
$ws.Range("H136").Value = 1091.3182
$ws.Range("I136").Value = 795.38464
$ws.Range("K136").Value = 2386.15392
$ws.Range("M136").Value = 163.8460800000003


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 131
$ws.Range("H131").Value = 3037.8474
$ws.Range("I131").Value = 521.6
$ws.Range("J131").Value = 3270.8333
$ws.Range("K131").Value = 1564.8
$ws.Range("L131").Value = 9812.499899999999
$ws.Range("M131").Value = 3475.2
$ws.Range("N131").Value = -19892.4999


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 74
$ws.Range("H74").Value = 20000
$ws.Range("J74").Value = 20000
$ws.Range("L74").Value = 20000
$ws.Range("N74").Value = -21996

# Row 77
$ws.Range("H77").Value = 20000
$ws.Range("J77").Value = 20000
$ws.Range("L77").Value = 60000
$ws.Range("N77").Value = -69984

# Row 82
$ws.Range("H82").Value = 27780250
$ws.Range("I82").Value = 50002000
$ws.Range("J82").Value = 3063.25
$ws.Range("K82").Value = 50002000
$ws.Range("L82").Value = 3063.25
$ws.Range("M82").Value = -50001639
$ws.Range("N82").Value = -3785.25

# Row 85
$ws.Range("H85").Value = 27780250
$ws.Range("I85").Value = 50002000
$ws.Range("J85").Value = 3063.25
$ws.Range("K85").Value = 50002000
$ws.Range("L85").Value = 3063.25
$ws.Range("M85").Value = -50000752
$ws.Range("N85").Value = -5559.25

# Row 136
$ws.Range("H136").Value = 12822875
$ws.Range("I136").Value = 2868
$ws.Range("K136").Value = 8604
$ws.Range("M136").Value = -6054


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 122
$ws.Range("H122").Value = 3249.625
$ws.Range("I122").Value = 3249.625
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9748.875
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7298.875
$ws.Range("N122").ClearContents()

# Row 136
$ws.Range("H136").Value = 3018.8774
$ws.Range("I136").Value = 2686.2974
$ws.Range("J136").Value = 4044.3333
$ws.Range("K136").Value = 8058.8922
$ws.Range("L136").Value = 12132.9999
$ws.Range("M136").Value = -5508.8922
$ws.Range("N136").Value = -17232.9999

